$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("79÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "63÷8=", 2) | Out-Null

$cell = $tbl.Cell(1, 2)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("69÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "90÷9=", 2) | Out-Null

$cell = $tbl.Cell(1, 3)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("30÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "97÷2=", 2) | Out-Null

$cell = $tbl.Cell(1, 4)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("46÷3=", $true, $false, $false, $false, $false, $true, 0, $false, "64÷3=", 2) | Out-Null

$cell = $tbl.Cell(1, 5)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("71÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "23÷3=", 2) | Out-Null

$cell = $tbl.Cell(5, 1)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("92÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "45÷5=", 2) | Out-Null

$cell = $tbl.Cell(5, 2)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("83÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "21÷9=", 2) | Out-Null

$cell = $tbl.Cell(5, 3)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("53÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "41÷2=", 2) | Out-Null

$cell = $tbl.Cell(5, 4)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("35÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "32÷5=", 2) | Out-Null

$cell = $tbl.Cell(5, 5)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("75÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "49÷9=", 2) | Out-Null

$cell = $tbl.Cell(9, 1)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("94÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "83÷7=", 2) | Out-Null

$cell = $tbl.Cell(9, 2)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("45÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "15÷2=", 2) | Out-Null

$cell = $tbl.Cell(9, 3)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("45÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "63÷8=", 2) | Out-Null

$cell = $tbl.Cell(9, 4)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("70÷2=", $true, $false, $false, $false, $false, $true, 0, $false, "14÷2=", 2) | Out-Null

$cell = $tbl.Cell(9, 5)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("72÷8=", $true, $false, $false, $false, $false, $true, 0, $false, "66÷6=", 2) | Out-Null

$cell = $tbl.Cell(13, 1)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("46÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "39÷2=", 2) | Out-Null

$cell = $tbl.Cell(13, 2)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("38÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "47÷8=", 2) | Out-Null

$cell = $tbl.Cell(13, 3)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("95÷6=", $true, $false, $false, $false, $false, $true, 0, $false, "69÷7=", 2) | Out-Null

$cell = $tbl.Cell(13, 4)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("51÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "48÷5=", 2) | Out-Null

$cell = $tbl.Cell(13, 5)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("95÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "66÷3=", 2) | Out-Null

$cell = $tbl.Cell(17, 1)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("92÷9=", $true, $false, $false, $false, $false, $true, 0, $false, "12÷2=", 2) | Out-Null

$cell = $tbl.Cell(17, 2)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("36÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "69÷4=", 2) | Out-Null

$cell = $tbl.Cell(17, 3)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("50÷4=", $true, $false, $false, $false, $false, $true, 0, $false, "63÷7=", 2) | Out-Null

$cell = $tbl.Cell(17, 4)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("64÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "53÷2=", 2) | Out-Null

$cell = $tbl.Cell(17, 5)
$cr = $cell.Range
$r = $d.Range($cr.Start, $cr.End)
$r.Find.Execute("94÷7=", $true, $false, $false, $false, $false, $true, 0, $false, "43÷5=", 2) | Out-Null
